$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($null -ne $val -and $val -is [string] -and $val -ne "") {
        $newVal = $val -replace "R2", "T2"
        $newVal = $newVal -replace "R1", "T"
        if ($newVal -ne $val) {
            $cell.Value2 = $newVal
        }
    }
}
